# Update sheet name and "through" date label from 07-03 to 07-04
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Name = "Through 2022-07-04"
$ws.Range("I1").Value = "2022 (through 07-04)"

# Update the July (row 8) and Total (row 14) values for the 2022 column (I)
$ws.Range("I8").Value = 27
$ws.Range("I14").Value = 833
